$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# The 9ad15719-6ce1-48cc-8569-036f14eacdc0 file has been handed back and is
# now in sync with en-US. Update the Overview sheet and the per-locale
# (zh-cn / de-de) detail sheets to reflect the handback.
# ---------------------------------------------------------------------------

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet: update status for the handed-back file ----------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusHandedBack
$overview.Range("C2").Value = $statusHandedBack

# --- zh-cn sheet -------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("B2").Value = $statusHandedBack

# Latest Target File (E2) + Latest Handback File (F2) are now populated,
# pointing at the same source/handoff targets as columns A2/C2.
$zhcn.Hyperlinks.Add(
    $zhcn.Range("E2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/1b714050bd19de299f2bc9b207fc0bdde29e62a9/e2e/9ad15719-6ce1-48cc-8569-036f14eacdc0.md",
    "",
    "",
    "9ad15719-6ce1-48cc-8569-036f14eacdc0.md"
) | Out-Null

$zhcn.Hyperlinks.Add(
    $zhcn.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5fe2827f2b618c4dbb3dd4f6b978994cc35b86d4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/9ad15719-6ce1-48cc-8569-036f14eacdc0.79be59d77f61d35d36ff06509bb7aaec296dacce.zh-cn.xlf",
    "",
    "",
    "9ad15719-6ce1-48cc-8569-036f14eacdc0.79be59d77f61d35d36ff06509bb7aaec296dacce.zh-cn.xlf"
) | Out-Null

# Latest Handback DateTime (G2)
$zhcn.Range("G2").Value = "2016-03-09 20:33:41"

# --- de-de sheet -------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("B2").Value = $statusHandedBack

$dede.Hyperlinks.Add(
    $dede.Range("E2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/1b714050bd19de299f2bc9b207fc0bdde29e62a9/e2e/9ad15719-6ce1-48cc-8569-036f14eacdc0.md",
    "",
    "",
    "9ad15719-6ce1-48cc-8569-036f14eacdc0.md"
) | Out-Null

$dede.Hyperlinks.Add(
    $dede.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/deaab2e6c5e4ca6e96d8f717b906980c1436fa92/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/9ad15719-6ce1-48cc-8569-036f14eacdc0.79be59d77f61d35d36ff06509bb7aaec296dacce.de-de.xlf",
    "",
    "",
    "9ad15719-6ce1-48cc-8569-036f14eacdc0.79be59d77f61d35d36ff06509bb7aaec296dacce.de-de.xlf"
) | Out-Null

# Latest Handback DateTime (G2)
$dede.Range("G2").Value = "2016-03-09 20:33:53"
